{"js": "// The document is a title paragraph (\"<date> <weekday>\") followed by a\n// single 20x5 table of arithmetic equations (\"a+b=c\" / \"a-b=c\"). The\n// commit bumps the date by a day and swaps in a new batch of equations,\n// cell-for-cell in row-major order (some old equation strings repeat, so\n// the mapping must be positional, not a global text search/replace).\n\nconst DATE_NEW = \"2023-10-16 Monday\";\n\n// Row-major replacement values for the 20x5 table (20 rows of 5 columns).\nconst TABLE_VALUES = [\n    [\"72-68=4\", \"6+45=51\", \"64-31=33\", \"31+64=95\", \"77-75=2\"],\n    [\"46+3=49\", \"27+66=93\", \"40+10=50\", \"11+57=68\", \"50+44=94\"],\n    [\"15+2=17\", \"60+36=96\", \"42+9=51\", \"34+39=73\", \"99-14=85\"],\n    [\"74-17=57\", \"70-60=10\", \"56-55=1\", \"84-43=41\", \"0+14=14\"],\n    [\"35+59=94\", \"71+9=80\", \"63-23=40\", \"79-4=75\", \"46-38=8\"],\n    [\"75-6=69\", \"86-73=13\", \"13+4=17\", \"32+47=79\", \"15+26=41\"],\n    [\"47+13=60\", \"83+0=83\", \"19+39=58\", \"94-78=16\", \"79-54=25\"],\n    [\"40-11=29\", \"59-6=53\", \"15-0=15\", \"53+13=66\", \"24+68=92\"],\n    [\"23+7=30\", \"20+47=67\", \"86-53=33\", \"37+41=78\", \"67-65=2\"],\n    [\"62-16=46\", \"73-57=16\", \"93-49=44\", \"99-60=39\", \"42-32=10\"],\n    [\"52-32=20\", \"67+23=90\", \"71-23=48\", \"10+22=32\", \"40+44=84\"],\n    [\"73-72=1\", \"73-31=42\", \"52-47=5\", \"5+52=57\", \"39+2=41\"],\n    [\"87+0=87\", \"0+17=17\", \"46-32=14\", \"6+29=35\", \"90-56=34\"],\n    [\"82-25=57\", \"81-40=41\", \"91+8=99\", \"71+10=81\", \"57-49=8\"],\n    [\"99-53=46\", \"65-27=38\", \"1+43=44\", \"11+37=48\", \"37+12=49\"],\n    [\"38-35=3\", \"85-67=18\", \"21+36=57\", \"78-46=32\", \"26+6=32\"],\n    [\"48+15=63\", \"36+25=61\", \"43-37=6\", \"52+26=78\", \"2+66=68\"],\n    [\"97-14=83\", \"22+3=25\", \"21+42=63\", \"75-51=24\", \"43+23=66\"],\n    [\"74+10=84\", \"31-2=29\", \"14+33=47\", \"12+49=61\", \"18+1=19\"],\n    [\"8+45=53\", \"6+84=90\", \"3+60=63\", \"72+22=94\", \"68+26=94\"],\n];\n\n// 1. Update the title paragraph's date/weekday text. insertText(..., replace)\n//    swaps the paragraph's text while keeping the existing run formatting\n//    (Arial, sz 30) of the first run.\nconst paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst titlePara = paras.items[0];\ntitlePara.insertText(DATE_NEW, Word.InsertLocation.replace);\n\n// 2. Update every cell of the single table in one shot, row-major, via the\n//    table's `values` property \u2014 this preserves per-cell run formatting\n//    (TimeNewRoman, sz 30) and cell/row/table structure.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = TABLE_VALUES;\n\nawait context.sync();\n", "ps1": "# The document is a title paragraph (\"<date> <weekday>\") followed by a\n# single 20x5 table of arithmetic equations (\"a+b=c\" / \"a-b=c\"). The\n# commit bumps the date by a day and swaps in a new batch of equations,\n# cell-for-cell in row-major order (some old equation strings repeat, so\n# the mapping must be positional, not a global Find/Replace).\n\n$dateNew = \"2023-10-16 Monday\"\n\n# Row-major replacement values for the 20x5 table (20 rows of 5 columns).\n$rows = @(\n    @(\"72-68=4\", \"6+45=51\", \"64-31=33\", \"31+64=95\", \"77-75=2\"),\n    @(\"46+3=49\", \"27+66=93\", \"40+10=50\", \"11+57=68\", \"50+44=94\"),\n    @(\"15+2=17\", \"60+36=96\", \"42+9=51\", \"34+39=73\", \"99-14=85\"),\n    @(\"74-17=57\", \"70-60=10\", \"56-55=1\", \"84-43=41\", \"0+14=14\"),\n    @(\"35+59=94\", \"71+9=80\", \"63-23=40\", \"79-4=75\", \"46-38=8\"),\n    @(\"75-6=69\", \"86-73=13\", \"13+4=17\", \"32+47=79\", \"15+26=41\"),\n    @(\"47+13=60\", \"83+0=83\", \"19+39=58\", \"94-78=16\", \"79-54=25\"),\n    @(\"40-11=29\", \"59-6=53\", \"15-0=15\", \"53+13=66\", \"24+68=92\"),\n    @(\"23+7=30\", \"20+47=67\", \"86-53=33\", \"37+41=78\", \"67-65=2\"),\n    @(\"62-16=46\", \"73-57=16\", \"93-49=44\", \"99-60=39\", \"42-32=10\"),\n    @(\"52-32=20\", \"67+23=90\", \"71-23=48\", \"10+22=32\", \"40+44=84\"),\n    @(\"73-72=1\", \"73-31=42\", \"52-47=5\", \"5+52=57\", \"39+2=41\"),\n    @(\"87+0=87\", \"0+17=17\", \"46-32=14\", \"6+29=35\", \"90-56=34\"),\n    @(\"82-25=57\", \"81-40=41\", \"91+8=99\", \"71+10=81\", \"57-49=8\"),\n    @(\"99-53=46\", \"65-27=38\", \"1+43=44\", \"11+37=48\", \"37+12=49\"),\n    @(\"38-35=3\", \"85-67=18\", \"21+36=57\", \"78-46=32\", \"26+6=32\"),\n    @(\"48+15=63\", \"36+25=61\", \"43-37=6\", \"52+26=78\", \"2+66=68\"),\n    @(\"97-14=83\", \"22+3=25\", \"21+42=63\", \"75-51=24\", \"43+23=66\"),\n    @(\"74+10=84\", \"31-2=29\", \"14+33=47\", \"12+49=61\", \"18+1=19\"),\n    @(\"8+45=53\", \"6+84=90\", \"3+60=63\", \"72+22=94\", \"68+26=94\")\n)\n\n$d = $word.ActiveDocument\n\n# 1. Update the title paragraph's date/weekday text. Setting Range.Text\n#    replaces the run's text while the run keeps its existing formatting\n#    (Arial, sz 30); the trailing paragraph mark is untouched.\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.Text = $dateNew\n\n# 2. Update every cell of the single table, row-major, one cell at a time.\n#    Cell.Range.Text replaces just the cell's content (the cell-end marker\n#    is handled automatically) and keeps the cell's run formatting\n#    (TimeNewRoman, sz 30) and the table's row/column structure intact.\n$t = $d.Tables.Item(1)\nfor ($r = 1; $r -le $rows.Count; $r++) {\n    $rowVals = $rows[$r - 1]\n    for ($c = 1; $c -le $rowVals.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $rowVals[$c - 1]\n    }\n}\n"}
